$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 204-205),
# pushing the existing rows 204-231 down to 206-233.
$ws.Rows("204:205").Insert()

# Row 204 (new): Región Metropolitana, week of 44491
$ws.Range("A204").Value = 6
$ws.Range("B204").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C204").Value = "Metropolitana"
$ws.Range("D204").Value = 44491
$ws.Range("E204").Value = 13
$ws.Range("F204").Value = 100112032
$ws.Range("G204").Value = "Zapallo italiano"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 400
$ws.Range("K204").Value = 14000
$ws.Range("L204").Value = 15000
$ws.Range("M204").Value = 14425
$ws.Range("N204").Value = "$/caja 50 unidades"
$ws.Range("O204").Value = "Región Metropolitana"
$ws.Range("P204").Value = 288
$ws.Range("Q204").Value = 50
$ws.Range("R204").Value = "Hortaliza"

# Row 205 (new): Región de O'Higgins, week of 44491
$ws.Range("A205").Value = 6
$ws.Range("B205").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C205").Value = "Metropolitana"
$ws.Range("D205").Value = 44491
$ws.Range("E205").Value = 13
$ws.Range("F205").Value = 100112032
$ws.Range("G205").Value = "Zapallo italiano"
$ws.Range("H205").Value = "Sin especificar"
$ws.Range("I205").Value = "Primera"
$ws.Range("J205").Value = 1000
$ws.Range("K205").Value = 14000
$ws.Range("L205").Value = 15000
$ws.Range("M205").Value = 14350
$ws.Range("N205").Value = "$/caja 50 unidades"
$ws.Range("O205").Value = "Región de O'Higgins"
$ws.Range("P205").Value = 287
$ws.Range("Q205").Value = 50
$ws.Range("R205").Value = "Hortaliza"
